$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

$values = @{
    2  = 102.33
    3  = 71.17
    4  = 65.43000000000001
    5  = 45.84
    6  = 44.57
    7  = 57.44
    8  = 81.7
    9  = 111.1
    10 = 109.76
    11 = 81.48999999999999
    12 = 26.73
    13 = 8.050000000000001
    14 = 3.85
    15 = 0.07000000000000001
    16 = 0.02
    17 = 1.39
    18 = 31.96
    19 = 88
    20 = 111.79
    21 = 128.94
    22 = 140
    23 = 142.18
    24 = 137.05
    25 = 121.38
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
